$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 522, shifting the
# existing rows 522:623 down to 524:625.
$ws.Range("A522:A523").EntireRow.Insert()

# New row 522 — Sandia, Extra, Región del Maule
$ws.Range("A522").Value = 10
$ws.Range("B522").Value = "Vega Modelo de Temuco"
$ws.Range("C522").Value = "La Araucanía"
$ws.Range("D522").Value = 44932
$ws.Range("E522").Value = 9
$ws.Range("F522").Value = 100112028
$ws.Range("G522").Value = "Sandia"
$ws.Range("H522").Value = "Sin especificar"
$ws.Range("I522").Value = "Extra"
$ws.Range("J522").Value = 4500
$ws.Range("K522").Value = 3500
$ws.Range("L522").Value = 3500
$ws.Range("M522").Value = 3500
$ws.Range("N522").Value = "$/unidad"
$ws.Range("O522").Value = "Región del Maule"
$ws.Range("P522").Value = 3500
$ws.Range("Q522").Value = 1
$ws.Range("R522").Value = "Hortaliza"

# New row 523 — Sandia, Primera, Región del Maule
$ws.Range("A523").Value = 10
$ws.Range("B523").Value = "Vega Modelo de Temuco"
$ws.Range("C523").Value = "La Araucanía"
$ws.Range("D523").Value = 44932
$ws.Range("E523").Value = 9
$ws.Range("F523").Value = 100112028
$ws.Range("G523").Value = "Sandia"
$ws.Range("H523").Value = "Sin especificar"
$ws.Range("I523").Value = "Primera"
$ws.Range("J523").Value = 6000
$ws.Range("K523").Value = 3000
$ws.Range("L523").Value = 3000
$ws.Range("M523").Value = 3000
$ws.Range("N523").Value = "$/unidad"
$ws.Range("O523").Value = "Región del Maule"
$ws.Range("P523").Value = 3000
$ws.Range("Q523").Value = 1
$ws.Range("R523").Value = "Hortaliza"

# Apply the date number format (style index 2 in the original workbook)
# to the two new date cells, matching the rest of column D.
$ws.Range("D522").NumberFormat = $ws.Range("D521").NumberFormat
$ws.Range("D523").NumberFormat = $ws.Range("D521").NumberFormat
